# New weekly data: insert 2 fresh rows at the top of the Agrícola del Norte
# S.A. de Arica - Maracuyá data block (currently rows 88-131), pushing the
# existing rows down by two positions (to 90-133), then populate the two
# newly-opened rows (88-89) with the latest week's prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 88-131 down to 90-133, opening up two blank rows at
# 88-89 (inheriting row 88's formatting, e.g. the date style on column D).
$ws.Rows("88:89").Insert()

# Row 88: Especial quality, 120 volume, new date + prices.
$ws.Range("A88").Value = 1
$ws.Range("B88").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C88").Value = "Arica y Parinacota"
$ws.Range("D88").Value = 44722
$ws.Range("E88").Value = 15
$ws.Range("F88").Value = "Fruta"
$ws.Range("G88").Value = 100108
$ws.Range("H88").Value = "Tropicales y subtropicales"
$ws.Range("I88").Value = 100108003
$ws.Range("J88").Value = "Maracuyá"
$ws.Range("K88").Value = "Sin especificar"
$ws.Range("L88").Value = "Especial"
$ws.Range("M88").Value = 120
$ws.Range("N88").Value = 21000
$ws.Range("O88").Value = 22000
$ws.Range("P88").Value = 21500
$ws.Range("Q88").Value = "$/caja 20 kilos"
$ws.Range("R88").Value = "Región de Arica y Parinacota"
$ws.Range("S88").Value = 1075
$ws.Range("T88").Value = 20

# Row 89: Primera quality, updated volume, new date + prices.
$ws.Range("A89").Value = 1
$ws.Range("B89").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C89").Value = "Arica y Parinacota"
$ws.Range("D89").Value = 44722
$ws.Range("E89").Value = 15
$ws.Range("F89").Value = "Fruta"
$ws.Range("G89").Value = 100108
$ws.Range("H89").Value = "Tropicales y subtropicales"
$ws.Range("I89").Value = 100108003
$ws.Range("J89").Value = "Maracuyá"
$ws.Range("K89").Value = "Sin especificar"
$ws.Range("L89").Value = "Primera"
$ws.Range("M89").Value = 140
$ws.Range("N89").Value = 18000
$ws.Range("O89").Value = 19000
$ws.Range("P89").Value = 18500
$ws.Range("Q89").Value = "$/caja 20 kilos"
$ws.Range("R89").Value = "Región de Arica y Parinacota"
$ws.Range("S89").Value = 925
$ws.Range("T89").Value = 20
